$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells that look like plain numbers are written
# back as text (matching the original inlineStr cell type) instead of being
# auto-converted to numeric values by Excel. We temporarily force a Text
# number format over the whole data range, write every new value, then
# restore the default "Normal" style so no stray formatting is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.431.18'
$ws.Range("D3").Value = '3.504.58'
$ws.Range("E3").Value = '  +4.10%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '586.10'
$ws.Range("E5").Value = '  +2.93%  '
$ws.Range("D6").Value = '147.64'
$ws.Range("E6").Value = '  +6.24%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +1.48%  '
$ws.Range("D9").Value = '7.71'
$ws.Range("E9").Value = '  +0.76%  '
$ws.Range("E10").Value = '  +4.57%  '
$ws.Range("D11").Value = '0.401'
$ws.Range("E11").Value = '  +5.33%  '
$ws.Range("D12").Value = '4.106.33'
$ws.Range("E12").Value = '  +4.18%  '
$ws.Range("D13").Value = '29.85'
$ws.Range("E13").Value = '  +7.93%  '
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("D15").Value = '3.498.55'
$ws.Range("E15").Value = '  +3.90%  '
$ws.Range("E16").Value = '  +4.60%  '
$ws.Range("D17").Value = '63.479.94'
$ws.Range("E17").Value = '  +4.19%  '
$ws.Range("E18").Value = '  +3.79%  '
$ws.Range("D19").Value = '14.30'
$ws.Range("E19").Value = '  +5.79%  '
$ws.Range("E20").Value = '  +6.94%  '
$ws.Range("D21").Value = '395.16'
$ws.Range("E21").Value = '  +3.86%  '
$ws.Range("D22").Value = '0.566'
$ws.Range("E22").Value = '  +3.44%  '
$ws.Range("D23").Value = '75.50'
$ws.Range("E23").Value = '  +0.31%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("E25").Value = '  +9.62%  '
$ws.Range("D26").Value = '3.646.19'
$ws.Range("E26").Value = '  +4.13%  '
$ws.Range("D27").Value = '0.189'
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("D28").Value = '7.82'
$ws.Range("E28").Value = '  +9.68%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.17%  '
$ws.Range("D30").Value = '8.28'
$ws.Range("E30").Value = '  +5.87%  '
$ws.Range("D31").Value = '2.16'
$ws.Range("E31").Value = '  +2.86%  '
$ws.Range("D32").Value = '1.42'
$ws.Range("E32").Value = '  +6.88%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").Value = '23.91'
$ws.Range("E34").Value = '  +4.34%  '
$ws.Range("E35").Value = '  +5.67%  '
$ws.Range("D36").Value = '32.59'
$ws.Range("E36").Value = '  +29.54%  '
$ws.Range("D37").Value = '5.36'
$ws.Range("E37").Value = '  +9.46%  '
$ws.Range("D38").Value = '172.41'
$ws.Range("E38").Value = '  +3.89%  '
$ws.Range("D39").Value = '1.58'
$ws.Range("E39").Value = '  +9.64%  '
$ws.Range("D40").Value = '3.542.58'
$ws.Range("E40").Value = '  +4.10%  '
$ws.Range("E41").Value = '  +1.94%  '
$ws.Range("D42").Value = '0.805'
$ws.Range("E42").Value = '  +4.16%  '
$ws.Range("E43").Value = '  +8.26%  '
$ws.Range("D44").Value = '4.52'
$ws.Range("E44").Value = '  +4.68%  '
$ws.Range("D45").Value = '42.63'
$ws.Range("E45").Value = '  +0.68%  '
$ws.Range("E46").Value = '  +10.49%  '
$ws.Range("D47").Value = '2.602.86'
$ws.Range("E47").Value = '  +6.18%  '
$ws.Range("D48").Value = '23.87'
$ws.Range("E48").Value = '  +7.90%  '
$ws.Range("D49").Value = '2.28'
$ws.Range("E49").Value = '  +12.90%  '
$ws.Range("D50").Value = '6.77'
$ws.Range("E50").Value = '  +2.61%  '
$ws.Range("D51").Value = '0.0271'
$ws.Range("E51").Value = '  +5.13%  '

# Restore default styling on the Price column now that the text is committed.
$ws.Range("D2:D51").Style = "Normal"
